$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; G = 5.488907176552729 }
    3 = @{ B = 1.505614041169197;  C = 0.3375848360084654; D = 0.7127328510149897; E = 6.48142807727062;   G = 9.037359805463273 }
    4 = @{ B = 0.7287194209349384; C = 1.65323645889881;   D = 3.082599426703578;  E = 0.4998867070740569; G = 5.964442013611383 }
    5 = @{ B = 0.7287194209349384; C = 0.3375848360084654; D = 0.7127328510149897; E = 6.48142807727062;   G = 8.260465185229014 }
    6 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 157.8057217802531;  E = 6.48142807727062;   G = 167.4460003575917 }
    7 = @{ B = 0.06328177979961902; C = 0.004309184025731883; D = 3.082599426703578; E = 0.4998867070740569; G = 3.650077097602987 }
    8 = @{ B = 0.7287194209349384; C = 0.3375848360084654; D = 0.7127328510149897; E = 6.48142807727062;   G = 8.260465185229014 }
    9 = @{ B = 1.505614041169197;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; G = 4.371470058157054 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
